$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 492, shifting the existing rows 492:512 down to 493:513
$ws.Rows("492:492").Insert()

# Populate the newly inserted row 492 with the new weekly price record
$ws.Cells.Item(492, 1).Value = 4
$ws.Cells.Item(492, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(492, 3).Value = "Los Lagos"
$ws.Cells.Item(492, 4).Value = "2023-08-09"
$ws.Cells.Item(492, 5).Value = 10
$ws.Cells.Item(492, 6).Value = 100112017
$ws.Cells.Item(492, 7).Value = "Apio"
$ws.Cells.Item(492, 8).Value = "Americana (o)"
$ws.Cells.Item(492, 9).Value = "Primera"
$ws.Cells.Item(492, 10).Value = 15
$ws.Cells.Item(492, 11).Value = 11000
$ws.Cells.Item(492, 12).Value = 11000
$ws.Cells.Item(492, 13).Value = 11000
$ws.Cells.Item(492, 14).Value = "`$/docena de matas"
$ws.Cells.Item(492, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(492, 16).Value = 1833
$ws.Cells.Item(492, 17).Value = 6
$ws.Cells.Item(492, 18).Value = "Hortaliza"
